$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 83892727
$ws.Range("B2").Value = 56411
$ws.Range("E2").Value = 100049
$ws.Range("F2").Value = "Spillkråka"
$ws.Range("G2").Value = "Dryocopus martius"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("M2").Value = "färska spår"
$ws.Range("P2").Value = "Stenmyrberget, Ång"
$ws.Range("Q2").Value = 715185.8176905914
$ws.Range("R2").Value = 7085533.152809725
$ws.Range("S2").Value = 5
$ws.Range("Y2").Value = "'2020-03-26"
$ws.Range("AA2").Value = "'2020-03-26"
$ws.Range("AW2").Value = "Emil Larsson"
$ws.Range("AX2").Value = "Emil Larsson"

$ws.Range("A3").Value = 80977
$ws.Range("B3").Value = 77540
$ws.Range("E3").Value = 185
$ws.Range("F3").Value = "Violettgrå tagellav"
$ws.Range("G3").Value = "Bryoria nadvornikiana"
$ws.Range("H3").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("M3").ClearContents()
$ws.Range("P3").Value = "Starrmyran, Ång"
$ws.Range("Q3").Value = 715363.1336768167
$ws.Range("R3").Value = 7085385.053305111
$ws.Range("S3").Value = 25
$ws.Range("Y3").Value = "'2007-02-13"
$ws.Range("AA3").Value = "'2007-02-13"
$ws.Range("AW3").Value = "Jonas F Grahn"
$ws.Range("AX3").Value = "Andreas Garpebring"
